$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 101-102; existing rows 101-123 shift down to 103-125.
$ws.Rows("101:102").Insert()

# New row 101
$ws.Range("A101").Value = 10
$ws.Range("B101").Value = "Vega Modelo de Temuco"
$ws.Range("C101").Value = "La Araucanía"
$ws.Range("D101").Value = 44468
$ws.Range("E101").Value = 9
$ws.Range("F101").Value = 100112013
$ws.Range("G101").Value = "Alcachofa"
$ws.Range("H101").Value = "Española"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 90
$ws.Range("K101").Value = 12000
$ws.Range("L101").Value = 13000
$ws.Range("M101").Value = 12556
$ws.Range("N101").Value = "$/caja 30 unidades"
$ws.Range("O101").Value = "Región Metropolitana"
$ws.Range("P101").Value = 419
$ws.Range("Q101").Value = 30
$ws.Range("R101").Value = "Hortaliza"

# New row 102
$ws.Range("A102").Value = 10
$ws.Range("B102").Value = "Vega Modelo de Temuco"
$ws.Range("C102").Value = "La Araucanía"
$ws.Range("D102").Value = 44468
$ws.Range("E102").Value = 9
$ws.Range("F102").Value = 100112013
$ws.Range("G102").Value = "Alcachofa"
$ws.Range("H102").Value = "Madrigal"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 80
$ws.Range("K102").Value = 12000
$ws.Range("L102").Value = 12000
$ws.Range("M102").Value = 12000
$ws.Range("N102").Value = "$/caja 40 unidades"
$ws.Range("O102").Value = "Región Metropolitana"
$ws.Range("P102").Value = 300
$ws.Range("Q102").Value = 40
$ws.Range("R102").Value = "Hortaliza"
